$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh A15's stored serial value (precision re-normalised by the
# upstream scraper when it re-read/re-wrote the sheet).
$ws.Range("A15").Value = 45815.39111230324

# Append the new day's scraped price row.
$ws.Range("A16").Value = 45816.39142756458
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C16").Value = "1Kg"
$ws.Range("D16").Value = "15,41€"
